$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values (row => new value) per commit "M10 Froze Encoder 12345"
$updates = @{
    2  = 6
    3  = 5
    4  = 7
    5  = 9
    6  = 6
    7  = 4
    8  = 7
    10 = 8
    12 = 4
    13 = 10
    14 = 4
    16 = 7
    18 = 10
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
